$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Malos Olores")

# Row 41: hyperlink text pointing at the elearning article (creates shared string for URL first)
$ws.Hyperlinks.Add($ws.Range("A41"), "https://elearning.industriallogic.com/gh/submit?Action=PageAction&album=recognizingSmells&path=recognizingSmells/featureEnvy/featureEnvyExample&devLanguage=Java")

# Row 40: "Referencias" heading
$ws.Range("A40").Value = "Referencias"

# Row 42: plain text URL (no hyperlink)
$ws.Range("A42").Value = "https://dzone.com/articles/code-smell-shot-surgery"

# Row 39: an empty styled cell
$ws.Range("A39").Value = $null

# Update the view: scroll so row 23 is at top, and select A42
$ws.Application.ActiveWindow.ScrollRow = 23
$ws.Range("A42").Select()
